$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column ("12-nov") right before the
#     existing "01-oct." column (currently column DQ / 121), shifting every
#     column from DQ..EU one position to the right (DR..EV). ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Inserting the entire column shifts DQ:EU -> DR:EV and keeps the header
# style (s="1") that was on the old DQ column.
$wsSpot.Range("DQ1").EntireColumn.Insert()

# New header cell (row 1) for the inserted column.
$wsSpot.Range("DQ1").Value = "12-nov"

# New data cells (rows 2-25) for the inserted column all show "-" (no data
# yet), matching every other not-yet-reached date column on the sheet.
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 121).Value = "-"
}

# --- Sheet "Gaz": append the next daily price row. ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(149, 1).NumberFormat = "@"
$wsGaz.Cells.Item(149, 1).Value = "2025-11-10"
$wsGaz.Cells.Item(149, 2).Value = 28.925

# --- Sheet "CO2": append the next daily price row. ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Cells.Item(149, 1).NumberFormat = "@"
$wsCO2.Cells.Item(149, 1).Value = "2025-11-10"
$wsCO2.Cells.Item(149, 2).Value = 79.88
